# Update "想去人数" (wanted-to-go count, column F) values on the
# "展览" and "全部类型" sheets to match the newly generated data snapshot.
# "全部类型" mirrors "展览"'s rows (plus one extra row), so the same
# per-event counts are updated on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1387
$ws1.Range("F5").Value = 109
$ws1.Range("F6").Value = 66
$ws1.Range("F7").Value = 11790
$ws1.Range("F8").Value = 4413
$ws1.Range("F10").Value = 44
$ws1.Range("F13").Value = 2556
$ws1.Range("F16").Value = 46
$ws1.Range("F17").Value = 5127
$ws1.Range("F19").Value = 187
$ws1.Range("F21").Value = 11364
$ws1.Range("F22").Value = 11313
$ws1.Range("F26").Value = 13
$ws1.Range("F28").Value = 22

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1387
$ws4.Range("F5").Value = 109
$ws4.Range("F6").Value = 66
$ws4.Range("F7").Value = 11790
$ws4.Range("F8").Value = 4413
$ws4.Range("F10").Value = 44
$ws4.Range("F13").Value = 2556
$ws4.Range("F17").Value = 46
$ws4.Range("F18").Value = 5127
$ws4.Range("F20").Value = 187
$ws4.Range("F22").Value = 11364
$ws4.Range("F23").Value = 11313
$ws4.Range("F27").Value = 13
$ws4.Range("F29").Value = 22
